# Scheduled-runner refresh of cached market-board figures (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ columns H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# sheets. Cells that have no corresponding market data (profit columns with no
# recipe-side price) are cleared rather than zeroed, matching upstream's removal of
# the <c> element for those rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 1000
$ws.Range("N4").Value = -1228
$ws.Range("M4").ClearContents()

$ws.Range("H88").Value = 2625
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2625
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2625
$ws.Range("N88").Value = -3437
$ws.Range("M88").ClearContents()

$ws.Range("H91").Value = 2625
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2625
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2625
$ws.Range("N91").Value = -5433
$ws.Range("M91").ClearContents()

$ws.Range("H111").Value = 9943
$ws.Range("I111").Value = 10029
$ws.Range("J111").Value = 9900
$ws.Range("K111").Value = 30087
$ws.Range("L111").Value = 29700
$ws.Range("M111").Value = -27020
$ws.Range("N111").Value = -35834

$ws.Range("H137").Value = 3520.8286
$ws.Range("I137").Value = 2455
$ws.Range("J137").Value = 7784.143
$ws.Range("K137").Value = 7365
$ws.Range("L137").Value = 23352.429
$ws.Range("M137").Value = -4815

$ws.Range("H138").Value = 2082
$ws.Range("I138").Value = 683
$ws.Range("J138").Value = 2272.7727
$ws.Range("K138").Value = 2049
$ws.Range("L138").Value = 6818.3181
$ws.Range("M138").Value = 3091
$ws.Range("N138").Value = -17098.3181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 334300.75
$ws.Range("I32").Value = 378742.62
$ws.Range("J32").Value = 15800.667
$ws.Range("K32").Value = 378742.62
$ws.Range("L32").Value = 15800.667
$ws.Range("M32").Value = -378455.62
$ws.Range("N32").Value = -16374.667

$ws.Range("H74").Value = 2868.1035
$ws.Range("I74").Value = 2610.8667
$ws.Range("J74").Value = 3143.7144
$ws.Range("K74").Value = 2610.8667
$ws.Range("L74").Value = 3143.7144
$ws.Range("M74").Value = -1736.8667
$ws.Range("N74").Value = -4891.7144

$ws.Range("H77").Value = 2868.1035
$ws.Range("I77").Value = 2610.8667
$ws.Range("J77").Value = 3143.7144
$ws.Range("K77").Value = 13054.3335
$ws.Range("L77").Value = 15718.572
$ws.Range("M77").Value = -8686.333500000001
$ws.Range("N77").Value = -24454.572

$ws.Range("H113").Value = 36532
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 36532
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 36532
$ws.Range("N113").Value = -45210
$ws.Range("M113").ClearContents()

$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1997.5
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1997.5
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1997.5
$ws.Range("N94").Value = -2899.5

$ws.Range("H134").Value = 2492.8286
$ws.Range("I134").Value = 2400.3333
$ws.Range("J134").Value = 2631.5715
$ws.Range("K134").Value = 7200.999899999999
$ws.Range("L134").Value = 7894.7145
$ws.Range("M134").Value = -4665.999899999999
$ws.Range("N134").Value = -12964.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 27636.545
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 27636.545
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 27636.545
$ws.Range("N4").Value = -27860.545

$ws.Range("H22").Value = 326
$ws.Range("I22").Value = 268.33334
$ws.Range("J22").Value = 412.5
$ws.Range("K22").Value = 268.33334
$ws.Range("L22").Value = 412.5
$ws.Range("M22").Value = 81.66665999999998
$ws.Range("N22").Value = -1112.5

$ws.Range("H31").Value = 7388.6924
$ws.Range("I31").Value = 2124.75
$ws.Range("J31").Value = 12929.685
$ws.Range("K31").Value = 2124.75
$ws.Range("L31").Value = 12929.685
$ws.Range("M31").Value = -1829.75
$ws.Range("N31").Value = -13519.685

$ws.Range("H34").Value = 7388.6924
$ws.Range("I34").Value = 2124.75
$ws.Range("J34").Value = 12929.685
$ws.Range("K34").Value = 2124.75
$ws.Range("L34").Value = 12929.685
$ws.Range("M34").Value = -1922.75
$ws.Range("N34").Value = -13333.685

$ws.Range("H111").Value = 24171.666
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 24171.666
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 24171.666
$ws.Range("N111").Value = -32351.666

$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws.Range("H132").Value = 11114144
$ws.Range("I132").Value = 3058.8572
$ws.Range("J132").Value = 20836344
$ws.Range("K132").Value = 9176.571599999999
$ws.Range("L132").Value = 62509032
$ws.Range("M132").Value = -6646.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 4644.25
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 4644.25
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 13932.75
$ws.Range("N94").Value = -15284.75
$ws.Range("M94").ClearContents()

$ws.Range("H119").Value = 3603
$ws.Range("I119").Value = 2904.5
$ws.Range("J119").Value = 5000
$ws.Range("K119").Value = 8713.5
$ws.Range("L119").Value = 15000
$ws.Range("M119").Value = -3875.5
$ws.Range("N119").Value = -24676

$ws.Range("H131").Value = 995.1458
$ws.Range("I131").Value = 450
$ws.Range("J131").Value = 1031.4889
$ws.Range("K131").Value = 1350
$ws.Range("L131").Value = 3094.4667
$ws.Range("M131").Value = 3690
$ws.Range("N131").Value = -13174.4667

$ws.Range("H140").Value = 1663.1818
$ws.Range("I140").Value = 1355.8
$ws.Range("J140").Value = 2321.8572
$ws.Range("K140").Value = 4067.4
$ws.Range("L140").Value = 6965.571599999999
$ws.Range("M140").Value = 1112.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5659.5947
$ws.Range("I70").Value = 5702.5835
$ws.Range("J70").Value = 5580.231
$ws.Range("K70").Value = 5702.5835
$ws.Range("L70").Value = 5580.231
$ws.Range("M70").Value = -5432.5835
$ws.Range("N70").Value = -6120.231

$ws.Range("H73").Value = 5659.5947
$ws.Range("I73").Value = 5702.5835
$ws.Range("J73").Value = 5580.231
$ws.Range("K73").Value = 5702.5835
$ws.Range("L73").Value = 5580.231
$ws.Range("M73").Value = -4766.5835
$ws.Range("N73").Value = -7452.231

$ws.Range("H97").Value = 3000
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 3000
$ws.Range("N97").Value = -3992
$ws.Range("M97").ClearContents()

$ws.Range("H122").Value = 4181.4
$ws.Range("I122").Value = 3503.5
$ws.Range("J122").Value = 4633.3335
$ws.Range("K122").Value = 10510.5
$ws.Range("L122").Value = 13900.0005
$ws.Range("M122").Value = -8060.5
$ws.Range("N122").Value = -18800.0005

$ws.Range("H140").Value = 58220
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 58220
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 58220
$ws.Range("N140").Value = -68580

$ws.Range("H141").Value = 40414.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 40414.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 40414.5
$ws.Range("N141").Value = -50774.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1714834.6
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 1714834.6
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 1714834.6
$ws.Range("N2").Value = -1715058.6

$ws.Range("H40").Value = 335468
$ws.Range("I40").Value = 1000004
$ws.Range("J40").Value = 3200
$ws.Range("K40").Value = 1000004
$ws.Range("L40").Value = 3200
$ws.Range("M40").Value = -999868
$ws.Range("N40").Value = -3472

$ws.Range("H93").Value = 14750.75
$ws.Range("I93").Value = 21201.2
$ws.Range("J93").Value = 4000
$ws.Range("K93").Value = 21201.2
$ws.Range("L93").Value = 4000
$ws.Range("M93").Value = -19953.2
$ws.Range("N93").Value = -6496

$ws.Range("H132").Value = 2116.7317
$ws.Range("I132").Value = 1210.5186
$ws.Range("J132").Value = 3864.4285
$ws.Range("K132").Value = 3631.5558
$ws.Range("L132").Value = 11593.2855
$ws.Range("M132").Value = -1101.5558
$ws.Range("N132").Value = -16653.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3755.318
$ws.Range("I136").Value = 3758.3572
$ws.Range("J136").Value = 3750
$ws.Range("K136").Value = 11275.0716
$ws.Range("L136").Value = 11250
$ws.Range("M136").Value = -8725.071599999999
$ws.Range("N136").Value = -16350
